$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "TSLA copy" to "TSLA"
$ws.Name = "TSLA"

# Append the new trading-day row (2020-12-31) after the existing last row (253)
$newRow = 254
$ws.Cells.Item($newRow, 1).Value = [DateTime]"2020-12-31"

# Open/High/Low/Close/AdjClose are stored as text in this sheet (like all the
# other rows), so force a Text number format before assigning, then drop the
# style back to Normal so the new cells don't carry a stray text format.
$priceCols = @(2, 3, 4, 5, 6)
$priceVals = @("699.989990", "718.719971", "691.119995", "705.669983", "705.669983")
for ($i = 0; $i -lt $priceCols.Length; $i++) {
    $cell = $ws.Cells.Item($newRow, $priceCols[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $priceVals[$i]
    $cell.Style = "Normal"
}

$ws.Cells.Item($newRow, 7).Value = 49649900
